$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Data changes
#
# Two new "building" rows (b10, b11) are inserted as rows 16-17, acting as
# windpark-type connection points on netNodes E1/E2 with a capacity formula
# of 400*1000. The pre-existing energyAsset rows (a1, a2, a3 - previously at
# rows 16-18, wired directly to the netNodes E2/E2/E3 with type "windpark")
# shift down to rows 18-20, get renumbered (ids 16,17,18) and are rewired so
# their "type2" is "windmolen" and their "parent" is the new building (b10 or
# b11) instead of the raw netNode.
# ---------------------------------------------------------------------------

# Make room: insert two blank rows above the existing energyAsset rows
# (old rows 16:17), which pushes the old rows 16,17,18 down to 18,19,20.
$ws.Rows("16:17").Insert()

# Seed the new "windmolen" lookup value first so it lands in the shared
# string table ahead of "b10"/"b11" (matches authoring order upstream).
$ws.Cells.Item(18, 5).Value = "windmolen"

# New row 16 - building "b10", type2 windpark, parent netNode E1
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "building"
$ws.Cells.Item(16, 3).Value = "b10"
$ws.Cells.Item(16, 4).Value = "windpark"
$ws.Cells.Item(16, 6).Value = "E1"
$ws.Cells.Item(16, 7).Formula = "=400*1000"
$ws.Cells.Item(16, 7).NumberFormat = "0.00E+00"

# New row 17 - building "b11", type2 windpark, parent netNode E2
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "building"
$ws.Cells.Item(17, 3).Value = "b11"
$ws.Cells.Item(17, 4).Value = "windpark"
$ws.Cells.Item(17, 6).Value = "E2"
$ws.Cells.Item(17, 7).Formula = "=400*1000"
$ws.Cells.Item(17, 7).NumberFormat = "0.00E+00"

# Row 18 (shifted from old row 16) - energyAsset a1, now parented to b10
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 3).Value = "a1"
$ws.Cells.Item(18, 6).Value = "b10"

# Row 19 (shifted from old row 17) - energyAsset a2, now parented to b10
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 3).Value = "a2"
$ws.Cells.Item(19, 5).Value = "windmolen"
$ws.Cells.Item(19, 6).Value = "b10"

# Row 20 (shifted from old row 18) - energyAsset a3, now parented to b11
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 3).Value = "a3"
$ws.Cells.Item(20, 5).Value = "windmolen"
$ws.Cells.Item(20, 6).Value = "b11"

# ---------------------------------------------------------------------------
# View / cosmetics
# ---------------------------------------------------------------------------

# New KPI column G gets a wider custom width
$ws.Columns.Item(7).ColumnWidth = 14.33203125

# Zoom in on the sheet and move the selection to the new KPI column
$wb.Windows.Item(1).Zoom = 130
$ws.Range("G5").Select()
